$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a logical (boolean) cell to A7: FALSE
$ws.Range("A7").Value = $false

# Move selection to A8 (next cell after the new data), matching the diff
$ws.Range("A8").Select()
